# The post "「今日は幸せでいることにしました」" (row 243) was removed from
# posts.xlsx. Delete that entire row so every following row shifts up by
# one, matching the new A1:C318 used range.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("243").Delete()
